$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting since values like "1.000"
# or "0.7124" would otherwise be auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.298.15'
$ws.Range('D3').Value = '1.871.95'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '0.7124'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').Value = '241.82'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '0.3110'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '0.07717'
$ws.Range('E9').Value = '  -3.33%  '
$ws.Range('D10').Value = '24.79'
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('D11').Value = '0.08401'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').Value = '1.868.58'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = '5.232'
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').Value = '0.7130'
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').Value = '91.31'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').Value = '29.310.66'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '5.949'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.000007987'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '243.91'
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('D20').Value = '2.122.44'
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('D21').Value = '13.17'
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '7.880'
$ws.Range('E23').Value = '  -2.19%  '
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '0.1646'
$ws.Range('E25').Value = '  +2.30%  '
$ws.Range('D26').Value = '163.85'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '9.001'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('E28').Value = '  +0.89%  '
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = '4.410'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('E31').Value = '  -3.82%  '
$ws.Range('D32').Value = '4.271'
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('D33').Value = '0.05165'
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('D34').Value = '0.7807'
$ws.Range('E34').Value = '  +7.75%  '
$ws.Range('D35').Value = '1.923'
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('D36').Value = '1.173'
$ws.Range('E36').Value = '  -2.38%  '
$ws.Range('D37').Value = '2.686'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '0.01862'
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('D39').Value = '2.711'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.164.43'
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '6.409'
$ws.Range('E41').Value = '  +4.21%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '73.34'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.8915'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '103.80'
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('D46').Value = '2.019.26'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '1.795'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.5183'
$ws.Range('E48').Value = '  -2.22%  '
$ws.Range('D49').Value = '9.422'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').Value = '0.00000000120'
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('D51').Value = '0.4309'
$ws.Range('E51').Value = '  -0.57%  '
